$wb = $excel.ActiveWorkbook

# ---- Sheet: Overview ----
$ws = $wb.Worksheets.Item('Overview')

# Set cell values (reordered rows after sorting by source file name)
$ws.Range('A1').Value = 'File Name'
$ws.Range('B1').Value = 'zh-cn'
$ws.Range('C1').Value = 'de-de'
$ws.Range('A2').Value = '0bcd2ce7-5fcc-48b5-94a0-47d35f7452df.md'
$ws.Range('B2').Value = 'In Translation'
$ws.Range('C2').Value = 'In Translation'
$ws.Range('A3').Value = 'b7d44263-9138-4e17-8083-9580d31ed741.md'
$ws.Range('B3').Value = 'In Translation'
$ws.Range('C3').Value = 'In Translation'
$ws.Range('A4').Value = 'da5e8a2c-8d37-416c-94e5-65f807b540ad.md'
$ws.Range('B4').Value = 'In Translation'
$ws.Range('C4').Value = 'In Translation'
$ws.Range('A5').Value = '.localization-config'
$ws.Range('B5').Value = 'Not to be localized'
$ws.Range('C5').Value = 'Not to be localized'

# Rebuild hyperlinks to match the new row order
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range('A2'), 'https://github.com/OpenLocalizationTest/oltest/blob/68b703065576b0821fc5d6d8974bfd03bfe02968/e2e/b7d44263-9138-4e17-8083-9580d31ed741.md', "", "", '0bcd2ce7-5fcc-48b5-94a0-47d35f7452df.md')
$ws.Hyperlinks.Add($ws.Range('A3'), 'https://github.com/OpenLocalizationTest/oltest/blob/68b703065576b0821fc5d6d8974bfd03bfe02968/e2e/da5e8a2c-8d37-416c-94e5-65f807b540ad.md', "", "", 'b7d44263-9138-4e17-8083-9580d31ed741.md')
$ws.Hyperlinks.Add($ws.Range('A4'), 'https://github.com/OpenLocalizationTest/oltest/blob/ca8097e91305999c2a5f9a892223f7602c984804/e2e/0bcd2ce7-5fcc-48b5-94a0-47d35f7452df.md', "", "", 'da5e8a2c-8d37-416c-94e5-65f807b540ad.md')
$ws.Hyperlinks.Add($ws.Range('A5'), 'https://github.com/OpenLocalizationTest/oltest/blob/ca8097e91305999c2a5f9a892223f7602c984804/.localization-config', "", "", '.localization-config')

# Restore the hyperlink font styling (cornflower-blue underline) that
# Hyperlinks.Add() does not reapply automatically on this engine
$ws.Range('A2').Font.Color = 15570276
$ws.Range('A2').Font.Underline = 2
$ws.Range('A3').Font.Color = 15570276
$ws.Range('A3').Font.Underline = 2
$ws.Range('A4').Font.Color = 15570276
$ws.Range('A4').Font.Underline = 2
$ws.Range('A5').Font.Color = 15570276
$ws.Range('A5').Font.Underline = 2

# ---- Sheet: zh-cn ----
$ws = $wb.Worksheets.Item('zh-cn')

# Set cell values (reordered rows after sorting by source file name)
$ws.Range('A1').Value = 'Source File Name'
$ws.Range('B1').Value = 'Status'
$ws.Range('C1').Value = 'Latest Handoff File'
$ws.Range('D1').Value = 'Latest Handoff Datetime'
$ws.Range('E1').Value = 'Latest Target File'
$ws.Range('F1').Value = 'Latest Handback File'
$ws.Range('G1').Value = 'Latest Handback DateTime'
$ws.Range('H1').Value = 'Handoff Reason'
$ws.Range('I1').Value = 'Dependency From'
$ws.Range('A2').Value = '0bcd2ce7-5fcc-48b5-94a0-47d35f7452df.md'
$ws.Range('B2').Value = 'In Translation'
$ws.Range('C2').Value = '0bcd2ce7-5fcc-48b5-94a0-47d35f7452df.a1ce4849b9565447c4a185bf00b08dbaf7048def.zh-cn.xlf'
$ws.Range('D2').Value = '2016-03-02 09:25:17'
$ws.Range('G2').Value = '0001-01-01 00:00:00'
$ws.Range('H2').Value = 'Include'
$ws.Range('A3').Value = 'b7d44263-9138-4e17-8083-9580d31ed741.md'
$ws.Range('B3').Value = 'In Translation'
$ws.Range('C3').Value = 'b7d44263-9138-4e17-8083-9580d31ed741.570d16d3f01576b9f1929dffc2b758d0d8b91ac2.zh-cn.xlf'
$ws.Range('D3').Value = '2016-03-02 09:23:41'
$ws.Range('G3').Value = '0001-01-01 00:00:00'
$ws.Range('H3').Value = 'Include'
$ws.Range('A4').Value = 'da5e8a2c-8d37-416c-94e5-65f807b540ad.md'
$ws.Range('B4').Value = 'In Translation'
$ws.Range('C4').Value = 'da5e8a2c-8d37-416c-94e5-65f807b540ad.09d3fd63c6b621b81bd9fbfbe8cf3948fa5fb65c.zh-cn.xlf'
$ws.Range('D4').Value = '2016-03-02 09:23:41'
$ws.Range('G4').Value = '0001-01-01 00:00:00'
$ws.Range('H4').Value = 'Include'
$ws.Range('A5').Value = '.localization-config'
$ws.Range('B5').Value = 'Not to be localized'
$ws.Range('D5').Value = '0001-01-01 00:00:00'
$ws.Range('G5').Value = '0001-01-01 00:00:00'
$ws.Range('H5').Value = 'Ignored'

# Rebuild hyperlinks to match the new row order
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range('A2'), 'https://github.com/OpenLocalizationTest/oltest/blob/68b703065576b0821fc5d6d8974bfd03bfe02968/e2e/b7d44263-9138-4e17-8083-9580d31ed741.md', "", "", '0bcd2ce7-5fcc-48b5-94a0-47d35f7452df.md')
$ws.Hyperlinks.Add($ws.Range('C2'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/abb3b05685450e9f50bf2a62cc43313e624cf6f4/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/b7d44263-9138-4e17-8083-9580d31ed741.570d16d3f01576b9f1929dffc2b758d0d8b91ac2.zh-cn.xlf', "", "", '0bcd2ce7-5fcc-48b5-94a0-47d35f7452df.a1ce4849b9565447c4a185bf00b08dbaf7048def.zh-cn.xlf')
$ws.Hyperlinks.Add($ws.Range('A3'), 'https://github.com/OpenLocalizationTest/oltest/blob/68b703065576b0821fc5d6d8974bfd03bfe02968/e2e/da5e8a2c-8d37-416c-94e5-65f807b540ad.md', "", "", 'b7d44263-9138-4e17-8083-9580d31ed741.md')
$ws.Hyperlinks.Add($ws.Range('C3'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/abb3b05685450e9f50bf2a62cc43313e624cf6f4/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/da5e8a2c-8d37-416c-94e5-65f807b540ad.09d3fd63c6b621b81bd9fbfbe8cf3948fa5fb65c.zh-cn.xlf', "", "", 'b7d44263-9138-4e17-8083-9580d31ed741.570d16d3f01576b9f1929dffc2b758d0d8b91ac2.zh-cn.xlf')
$ws.Hyperlinks.Add($ws.Range('A4'), 'https://github.com/OpenLocalizationTest/oltest/blob/ca8097e91305999c2a5f9a892223f7602c984804/e2e/0bcd2ce7-5fcc-48b5-94a0-47d35f7452df.md', "", "", 'da5e8a2c-8d37-416c-94e5-65f807b540ad.md')
$ws.Hyperlinks.Add($ws.Range('C4'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9ae40d5678e013b399956086defa2bd1524beb12/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/0bcd2ce7-5fcc-48b5-94a0-47d35f7452df.a1ce4849b9565447c4a185bf00b08dbaf7048def.zh-cn.xlf', "", "", 'da5e8a2c-8d37-416c-94e5-65f807b540ad.09d3fd63c6b621b81bd9fbfbe8cf3948fa5fb65c.zh-cn.xlf')
$ws.Hyperlinks.Add($ws.Range('A5'), 'https://github.com/OpenLocalizationTest/oltest/blob/ca8097e91305999c2a5f9a892223f7602c984804/.localization-config', "", "", '.localization-config')

# Restore the hyperlink font styling (cornflower-blue underline) that
# Hyperlinks.Add() does not reapply automatically on this engine
$ws.Range('A2').Font.Color = 15570276
$ws.Range('A2').Font.Underline = 2
$ws.Range('C2').Font.Color = 15570276
$ws.Range('C2').Font.Underline = 2
$ws.Range('A3').Font.Color = 15570276
$ws.Range('A3').Font.Underline = 2
$ws.Range('C3').Font.Color = 15570276
$ws.Range('C3').Font.Underline = 2
$ws.Range('A4').Font.Color = 15570276
$ws.Range('A4').Font.Underline = 2
$ws.Range('C4').Font.Color = 15570276
$ws.Range('C4').Font.Underline = 2
$ws.Range('A5').Font.Color = 15570276
$ws.Range('A5').Font.Underline = 2

# ---- Sheet: de-de ----
$ws = $wb.Worksheets.Item('de-de')

# Set cell values (reordered rows after sorting by source file name)
$ws.Range('A1').Value = 'Source File Name'
$ws.Range('B1').Value = 'Status'
$ws.Range('C1').Value = 'Latest Handoff File'
$ws.Range('D1').Value = 'Latest Handoff Datetime'
$ws.Range('E1').Value = 'Latest Target File'
$ws.Range('F1').Value = 'Latest Handback File'
$ws.Range('G1').Value = 'Latest Handback DateTime'
$ws.Range('H1').Value = 'Handoff Reason'
$ws.Range('I1').Value = 'Dependency From'
$ws.Range('A2').Value = '0bcd2ce7-5fcc-48b5-94a0-47d35f7452df.md'
$ws.Range('B2').Value = 'In Translation'
$ws.Range('C2').Value = '0bcd2ce7-5fcc-48b5-94a0-47d35f7452df.a1ce4849b9565447c4a185bf00b08dbaf7048def.de-de.xlf'
$ws.Range('D2').Value = '2016-03-02 09:25:29'
$ws.Range('G2').Value = '0001-01-01 00:00:00'
$ws.Range('H2').Value = 'Include'
$ws.Range('A3').Value = 'b7d44263-9138-4e17-8083-9580d31ed741.md'
$ws.Range('B3').Value = 'In Translation'
$ws.Range('C3').Value = 'b7d44263-9138-4e17-8083-9580d31ed741.570d16d3f01576b9f1929dffc2b758d0d8b91ac2.de-de.xlf'
$ws.Range('D3').Value = '2016-03-02 09:23:57'
$ws.Range('G3').Value = '0001-01-01 00:00:00'
$ws.Range('H3').Value = 'Include'
$ws.Range('A4').Value = 'da5e8a2c-8d37-416c-94e5-65f807b540ad.md'
$ws.Range('B4').Value = 'In Translation'
$ws.Range('C4').Value = 'da5e8a2c-8d37-416c-94e5-65f807b540ad.09d3fd63c6b621b81bd9fbfbe8cf3948fa5fb65c.de-de.xlf'
$ws.Range('D4').Value = '2016-03-02 09:23:57'
$ws.Range('G4').Value = '0001-01-01 00:00:00'
$ws.Range('H4').Value = 'Include'
$ws.Range('A5').Value = '.localization-config'
$ws.Range('B5').Value = 'Not to be localized'
$ws.Range('D5').Value = '0001-01-01 00:00:00'
$ws.Range('G5').Value = '0001-01-01 00:00:00'
$ws.Range('H5').Value = 'Ignored'

# Rebuild hyperlinks to match the new row order
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range('A2'), 'https://github.com/OpenLocalizationTest/oltest/blob/68b703065576b0821fc5d6d8974bfd03bfe02968/e2e/b7d44263-9138-4e17-8083-9580d31ed741.md', "", "", '0bcd2ce7-5fcc-48b5-94a0-47d35f7452df.md')
$ws.Hyperlinks.Add($ws.Range('C2'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dd8ddf4325ab0f5e12e25a4e46a9e1fd82e30a8e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/b7d44263-9138-4e17-8083-9580d31ed741.570d16d3f01576b9f1929dffc2b758d0d8b91ac2.de-de.xlf', "", "", '0bcd2ce7-5fcc-48b5-94a0-47d35f7452df.a1ce4849b9565447c4a185bf00b08dbaf7048def.de-de.xlf')
$ws.Hyperlinks.Add($ws.Range('A3'), 'https://github.com/OpenLocalizationTest/oltest/blob/68b703065576b0821fc5d6d8974bfd03bfe02968/e2e/da5e8a2c-8d37-416c-94e5-65f807b540ad.md', "", "", 'b7d44263-9138-4e17-8083-9580d31ed741.md')
$ws.Hyperlinks.Add($ws.Range('C3'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dd8ddf4325ab0f5e12e25a4e46a9e1fd82e30a8e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/da5e8a2c-8d37-416c-94e5-65f807b540ad.09d3fd63c6b621b81bd9fbfbe8cf3948fa5fb65c.de-de.xlf', "", "", 'b7d44263-9138-4e17-8083-9580d31ed741.570d16d3f01576b9f1929dffc2b758d0d8b91ac2.de-de.xlf')
$ws.Hyperlinks.Add($ws.Range('A4'), 'https://github.com/OpenLocalizationTest/oltest/blob/ca8097e91305999c2a5f9a892223f7602c984804/e2e/0bcd2ce7-5fcc-48b5-94a0-47d35f7452df.md', "", "", 'da5e8a2c-8d37-416c-94e5-65f807b540ad.md')
$ws.Hyperlinks.Add($ws.Range('C4'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6e90082ca96724b00258311f574905d24abcd7ef/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/0bcd2ce7-5fcc-48b5-94a0-47d35f7452df.a1ce4849b9565447c4a185bf00b08dbaf7048def.de-de.xlf', "", "", 'da5e8a2c-8d37-416c-94e5-65f807b540ad.09d3fd63c6b621b81bd9fbfbe8cf3948fa5fb65c.de-de.xlf')
$ws.Hyperlinks.Add($ws.Range('A5'), 'https://github.com/OpenLocalizationTest/oltest/blob/ca8097e91305999c2a5f9a892223f7602c984804/.localization-config', "", "", '.localization-config')

# Restore the hyperlink font styling (cornflower-blue underline) that
# Hyperlinks.Add() does not reapply automatically on this engine
$ws.Range('A2').Font.Color = 15570276
$ws.Range('A2').Font.Underline = 2
$ws.Range('C2').Font.Color = 15570276
$ws.Range('C2').Font.Underline = 2
$ws.Range('A3').Font.Color = 15570276
$ws.Range('A3').Font.Underline = 2
$ws.Range('C3').Font.Color = 15570276
$ws.Range('C3').Font.Underline = 2
$ws.Range('A4').Font.Color = 15570276
$ws.Range('A4').Font.Underline = 2
$ws.Range('C4').Font.Color = 15570276
$ws.Range('C4').Font.Underline = 2
$ws.Range('A5').Font.Color = 15570276
$ws.Range('A5').Font.Underline = 2
